$d = $word.ActiveDocument

# Each question's five answer choices (A-E) occupy five consecutive
# paragraphs in the "Questions" section. We address each answer choice
# by its absolute paragraph index so we never risk touching the wrong
# occurrence of a repeated "X: Correct/Incorrect Answer" string (there
# are many duplicates of e.g. "D: Correct Answer" across questions).
#
# For each entry below: OldText is what the paragraph currently holds
# (used as a sanity check) and NewText is what it should become.

$changes = @(
    @{ Index = 12; OldText = "D: Correct Answer";   NewText = "D: Incorrect Answer" },  # Q1 D
    @{ Index = 13; OldText = "E: Incorrect Answer"; NewText = "E: Correct Answer" },    # Q1 E
    @{ Index = 22; OldText = "B: Incorrect Answer"; NewText = "B: Correct Answer" },    # Q3 B
    @{ Index = 24; OldText = "D: Correct Answer";   NewText = "D: Incorrect Answer" },  # Q3 D
    @{ Index = 29; OldText = "C: Correct Answer";   NewText = "C: Incorrect Answer" },  # Q4 C
    @{ Index = 31; OldText = "E: Incorrect Answer"; NewText = "E: Correct Answer" },    # Q4 E
    @{ Index = 34; OldText = "B: Correct Answer";   NewText = "B: Incorrect Answer" },  # Q5 B
    @{ Index = 36; OldText = "D: Incorrect Answer"; NewText = "D: Correct Answer" },    # Q5 D
    @{ Index = 41; OldText = "C: Incorrect Answer"; NewText = "C: Correct Answer" },    # Q6 C
    @{ Index = 42; OldText = "D: Correct Answer";   NewText = "D: Incorrect Answer" },  # Q6 D
    @{ Index = 45; OldText = "A: Correct Answer";   NewText = "A: Incorrect Answer" },  # Q7 A
    @{ Index = 48; OldText = "D: Incorrect Answer"; NewText = "D: Correct Answer" },    # Q7 D
    @{ Index = 51; OldText = "A: Incorrect Answer"; NewText = "A: Correct Answer" },    # Q8 A
    @{ Index = 54; OldText = "D: Correct Answer";   NewText = "D: Incorrect Answer" },  # Q8 D
    @{ Index = 60; OldText = "D: Incorrect Answer"; NewText = "D: Correct Answer" },    # Q9 D
    @{ Index = 61; OldText = "E: Correct Answer";   NewText = "E: Incorrect Answer" },  # Q9 E
    @{ Index = 65; OldText = "C: Correct Answer";   NewText = "C: Incorrect Answer" },  # Q10 C
    @{ Index = 66; OldText = "D: Incorrect Answer"; NewText = "D: Correct Answer" }     # Q10 D
)

foreach ($c in $changes) {
    $p = $d.Paragraphs.Item($c.Index)
    $current = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $c.OldText) {
        throw "Paragraph $($c.Index) expected '$($c.OldText)' but found '$current'"
    }
    $p.Range.Text = $c.NewText
}
